$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with fresh crypto data.
# Numeric-looking Price strings are forced back to text so the cell keeps
# its original string representation (e.g. "214.40" instead of 214.4).

$ws.Range("D2").Value = "28.431.12"
$ws.Range("E2").Value = "  +4.66%  "

$ws.Range("D3").Value = "1.588.91"
$ws.Range("E3").Value = "  +1.22%  "

$ws.Range("E4").Value = "  -0.54%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "214.40"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +1.52%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.498"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +1.09%  "

$ws.Range("E7").Value = "  -0.71%  "

$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "23.92"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  +8.22%  "

$ws.Range("E9").Value = "  +0.73%  "

$ws.Range("E10").Value = "  +0.25%  "

$ws.Range("E11").Value = "  +1.75%  "

$ws.Range("D12").Value = "1.815.23"
$ws.Range("E12").Value = "  +1.25%  "

$ws.Range("D13").Value = "1.589.04"
$ws.Range("E13").Value = "  +1.20%  "

$ws.Range("E14").Value = "  +0.10%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "0.533"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  +2.79%  "

$ws.Range("D16").Value = "28.434.16"
$ws.Range("E16").Value = "  +4.28%  "

$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "63.11"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +1.40%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "232.86"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +7.07%  "

$ws.Range("E19").Value = "  +0.64%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "7.50"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  -0.04%  "

$ws.Range("E21").Value = "  -0.55%  "

$ws.Range("E22").Value = "  -0.80%  "

$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "9.42"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +1.97%  "

$ws.Range("E24").Value = "  +0.30%  "

$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "152.06"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  -0.99%  "

$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "15.28"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +1.39%  "

$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "6.63"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.35%  "

$ws.Range("E28").Value = "  +0.72%  "

$ws.Range("E29").Value = "  -0.54%  "

$ws.Range("E30").Value = "  +0.05%  "

$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.0474"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.48%  "

$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.25"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.07%  "

$ws.Range("E33").Value = "  -0.14%  "

$ws.Range("D34").Value = "1.416.44"
$ws.Range("E34").Value = "  -1.97%  "

$ws.Range("E35").Value = "  -1.28%  "

$ws.Range("E36").Value = "  -5.32%  "

$ws.Range("E37").Value = "  -0.52%  "

$ws.Range("E38").Value = "  +0.39%  "

$ws.Range("E39").Value = "  +7.93%  "

$ws.Range("E40").Value = "  +1.87%  "

$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "0.819"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  +1.07%  "

$ws.Range("E42").Value = "  -2.68%  "

$ws.Range("E43").Value = "  -0.76%  "

$ws.Range("E44").Value = "  -2.47%  "

$ws.Range("E45").Value = "  +5.69%  "

$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "64.63"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.06%  "

$ws.Range("D47").Value = "1.728.42"
$ws.Range("E47").Value = "  +1.41%  "

$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "87.53"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  +1.71%  "

$ws.Range("E49").Value = "  +5.76%  "

$ws.Range("E50").Value = "  -0.73%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "39.39"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +16.17%  "
